# Auto-committed on 2023/08/04 週五 16:06:39.95
#
# Updates the PfItDetail.xlsx DB-layout sheets:
#  - DBD!G35 remark gets a "(不產媒體)" suffix
#  - DBS key-id row gets renamed (findBormNoEq -> findBormNoFirst) and its
#    order-by strings gain the LogNo tie-breaker
#  - selection / active sheet ends up parked on DBS!C11

$wb = $excel.ActiveWorkbook

$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# --- DBD sheet: amend the remark for MediaFg (row 35, column G) ---
$wsDBD.Range("G35").Value = "3.保費檢核結果為Y時已追回撥款，還款不用追回(不產媒體)"

# --- DBS sheet: rename the key id and extend the two ORDER BY clauses ---
# (order of first-use matters for shared-string table layout, so write
#  C11 before C10 to match the canonical append order)
$wsDBS.Range("A11").Value = "findBormNoFirst"
$wsDBS.Range("C11").Value = "PerfDate ASC, LogNo ASC"
$wsDBS.Range("C10").Value = "PerfDate Desc , LogNo DESC"

# --- Selection / active sheet bookkeeping ---
$wsDBD.Select() | Out-Null
$wsDBD.Range("B9").Select() | Out-Null

$wsDBS.Select() | Out-Null
$wsDBS.Range("C11").Select() | Out-Null
